# Add a new "2022-Q3" sheet right after "总计" (shifting 2022-Q2 .. 2021-Q2 one
# slot to the right, which happens automatically since they keep their own
# names), and insert the corresponding new summary row on "总计".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new quarter sheet right after "总计".
# ---------------------------------------------------------------------------
$zongji = $wb.Worksheets.Item(1)
$q3 = $wb.Worksheets.Add($null, $zongji)
$q3.Name = "2022-Q3"

# Header row (same headers used by the other quarterly sheets).
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"
$q3.Range("B1:H1").Style = "Normal"
$q3.Range("B1:H1").HorizontalAlignment = -4108
$q3.Range("B1:H1").VerticalAlignment = -4160
$q3.Range("B1:H1").Font.Bold = $true
$q3.Range("B1:H1").Borders.LineStyle = 1

# Data rows. Columns B-G are stored as text (fund codes / numbers with fixed
# formatting), so force the text number format before writing, then drop back
# to the default style so no stray "@" format sticks around.
$q3TextCols = $q3.Range("B2:G4")
$q3TextCols.NumberFormat = "@"

$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "005702"
$q3.Range("C2").Value = "恒生前海港股通高股息低波动指数"
$q3.Range("D2").Value = "0.20"
$q3.Range("E2").Value = "94.22"
$q3.Range("F2").Value = "3.44"
$q3.Range("G2").Value = "0.0069"
$q3.Range("H2").Value = 2

$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "501307"
$q3.Range("C3").Value = "银河中证沪港深高股息指数（LOF）A"
$q3.Range("D3").Value = "0.15"
$q3.Range("E3").Value = "90.33"
$q3.Range("F3").Value = "1.28"
$q3.Range("G3").Value = "0.0019"
$q3.Range("H3").Value = 8

$q3.Range("A4").Value = 2
$q3.Range("B4").Value = "501308"
$q3.Range("C4").Value = "银河中证沪港深高股息指数（LOF）C"
$q3.Range("D4").Value = "0.01"
$q3.Range("E4").Value = "90.33"
$q3.Range("F4").Value = "1.28"
$q3.Range("G4").Value = "0.0001"
$q3.Range("H4").Value = 8

$q3TextCols.Style = "Normal"

$q3.Range("A2:A4").Style = "Normal"
$q3.Range("A2:A4").HorizontalAlignment = -4108
$q3.Range("A2:A4").VerticalAlignment = -4160
$q3.Range("A2:A4").Font.Bold = $true
$q3.Range("A2:A4").Borders.LineStyle = 1

$q3.Range("A1").Select()

# ---------------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert a new row for 2022-Q3 above the
#    existing 2022-Q2 row, pushing everything else down by one.
# ---------------------------------------------------------------------------
$zongji.Rows.Item(2).Insert()

$zongji.Range("A2").Value = 0
$zongji.Range("A2").HorizontalAlignment = -4108
$zongji.Range("A2").VerticalAlignment = -4160
$zongji.Range("A2").Font.Bold = $true
$zongji.Range("A2").Borders.LineStyle = 1

$zongji.Range("B2").Value = "2022-Q3"
$zongji.Range("C2").Value = 3
$zongji.Range("D2").Value = 0.01

$zongji.Range("A1").Select()
